# "Practiced IELTS 11 Test 3 and 4" - fill in the IELTS11_Test4 results row
# and tidy up a few leftover cell-format inconsistencies further down the
# scoring table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# Row 16 (IELTS11_Test4): add the test date and the scores that were
# entered for this sitting. Copy the date format from the row above so
# the new cell reuses the existing "short date" style instead of Excel
# minting a brand-new (duplicate) number format.
# ---------------------------------------------------------------------
$ws.Range("C15").Copy()
$ws.Range("C16").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("C16").Value = 45452

$ws.Range("E16").Value = 30
$ws.Range("F16").Formula = '=IFERROR(INDEX(Sheet2!$F$5:$F$20, MATCH(Table1[[#This Row],[Lis_Mark]], Sheet2!$D$5:$D$20, 1)),"No Grade")'

$ws.Range("G16").Value = 26
$ws.Range("H16").Formula = '=IFERROR(INDEX(Sheet2!$F$5:$F$20, MATCH(Table1[[#This Row],[Read_Mark]], Sheet2!$D$5:$D$20, 1)),"No Grade")'

$ws.Range("I16").Value = 1.1
# K16 already holds the shared "(F+H+I+J)/3" formula for this row, so it
# recalculates to the new Overall score automatically.

# ---------------------------------------------------------------------
# H22:H25 had drifted onto a near-duplicate style (same font, just with a
# redundant applyNumberFormat flag). Re-copy the format that the rest of
# the column (H17/H18/H20/H21) already uses so they share the same style
# again instead of a stray almost-identical one.
# ---------------------------------------------------------------------
$ws.Range("H21").Copy()
$ws.Range("H22:H25").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# Leave the cursor where it was left in the saved workbook.
$ws.Range("J13").Select()
